# Updated symbol list on Fri Feb 10 17:36:55 UTC 2023 with GitHub Actions
# Refreshes Price/Volume(1h) figures for existing coins, and shifts the
# UpBots row up to position 26 (with refreshed data) while the other rows
# that used to sit between Spectre.aiUtilityToken and UpBots each move down
# by one slot, pushing BNIXToken into row 34 with placeholder "--"/"--%"
# values.
#
# Values are assigned with a leading "'" to force Excel to store them as
# text (matching the original inlineStr/text cell type instead of letting
# COM auto-convert numeric- or percent-looking strings into real numbers),
# and Style is reset to "Normal" right afterwards so no stray
# "quote prefix" number format is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-3.16%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'40.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-3.28%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.039"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-2.71%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07589"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-6.43%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.247"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-2.72%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.591"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-8.87%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9041"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-2.82%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.09941"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-11.60%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1763"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-5.24%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09006"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-3.09%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04394"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-4.02%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1053"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.14%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001256"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-1.94%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005826"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.40%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.367"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.45%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'-3.96%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'-2.79%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.844"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-7.67%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1352"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-2.16%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'9.75%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'-0.63%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'-2.11%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004066"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001303"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'6.32%"
$ws.Range("E25").Style = "Normal"
$ws.Range("B26").Value = "'UpBots"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'0.0003017"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'1.08%"
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "'Spectre.aiUtilityToken"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut"
$ws.Range("C27").Style = "Normal"
$ws.Range("B28").Value = "'LegolasExchange"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo"
$ws.Range("C28").Style = "Normal"
$ws.Range("B29").Value = "'BitZToken"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz"
$ws.Range("C29").Style = "Normal"
$ws.Range("B30").Value = "'Birake"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir"
$ws.Range("C30").Style = "Normal"
$ws.Range("B31").Value = "'NashExchange"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex"
$ws.Range("C31").Style = "Normal"
$ws.Range("B32").Value = "'AAXToken"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab"
$ws.Range("C32").Style = "Normal"
$ws.Range("B33").Value = "'CenX"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/V4XJUvLQb+cenx-cenx"
$ws.Range("C33").Style = "Normal"
$ws.Range("B34").Value = "'BNIXToken"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'--"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'--%"
$ws.Range("E34").Style = "Normal"
$ws.Range("D38").Value = "'0.02392"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-7.96%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05125"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-6.56%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007853"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-2.58%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1302"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-6.53%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007096"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'8.95%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.001953"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-6.72%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008383"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'1.46%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3319"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-3.80%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006460"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-3.91%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E48").Value = "'-26.88%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.005713"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'68.01%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D51").Value = "'0.0002006"
$ws.Range("D51").Style = "Normal"
